# IMLGS-223 Sample Comments and Sample Lake add to spreadsheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Lake" column (AR) ---
$ws.Range("AR1").Value = "Lake"
$ws.Range("AR2").Value = "Lake 4"
$ws.Range("AR3").Value = "Lake 4"
$ws.Range("AR4").Value = "Lake 1"
$ws.Range("AR5").Value = "Lake 2"
$ws.Range("AR6").Value = "Lake 3"

# --- New "Sample Comments" column (AS) ---
$ws.Range("AS1").Value = "Sample Comments"
$ws.Range("AS2").Value = "Sample comment 4"
$ws.Range("AS3").Value = "Sample comment 4"
$ws.Range("AS6").Value = "Sample comment 3"
$ws.Range("AS5").Value = "Sample comment 2"
$ws.Range("AS4").Value = "Sample comment 1"

# Data rows for AR already inherit the existing centered/Arial style from the
# previously-blank AR2:AR6 cells. Apply the same formatting to the new AS
# data cells by copying the format from the AR column.
$ws.Range("AR2:AR6").Copy()
$ws.Range("AS2:AS6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column widths ---
$ws.Range("B:B").ColumnWidth = 18.666666666666668
$ws.Range("AS:AS").ColumnWidth = 14.333333333333334

# --- Selection / view state ---
$ws.Range("AS5").Select()
